$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 148
$ws.Range("I2").Value = 408
$ws.Range("J2").Value = 1591
$ws.Range("K2").Value = 12
$ws.Range("L2").Value = 417
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 258
$ws.Range("R2").Value = 30
$ws.Range("S2").Value = 167
$ws.Range("T2").Value = 239
$ws.Range("U2").Value = 16
$ws.Range("V2").Value = 2452
$ws.Range("W2").Value = 0
$ws.Range("Y2").Value = 2
$ws.Range("AA2").Value = 19
